$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above the old footer row (row 9) so that it becomes
#    row 10. Excel's default Insert() copies the formatting of the row above
#    (row 8) into the freshly inserted row 9, including the P8:Q8 merge,
#    which is exactly the merge we need for the new "total" cell P9:Q9 -
#    we'll just need to fix its height and give it a value.
# ---------------------------------------------------------------------------
$ws.Rows(9).Insert()

# ---------------------------------------------------------------------------
# 2. Row 7 - fill in the (previously empty) item-detail cells.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "بخاخ ريد الكبير "
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7:M7").Value = "0"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "100.00"
$ws.Range("P7").Value = "100.0000"
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# ---------------------------------------------------------------------------
# 3. Row 8 - turn it into a second item-detail row, matching row 7's layout.
#    First undo the inherited P8:Q8 merge (that belonged to the old footer
#    total cell, which has now moved down to row 9), then build the same
#    merge/style layout row 7 has.
# ---------------------------------------------------------------------------
$ws.Range("P8:Q8").UnMerge()

$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A8").Value = 2
$ws.Range("C8:G8").NumberFormat = "@"
$ws.Range("C8").Value = "سائل ريد"
$ws.Range("H8:K8").NumberFormat = "@"
$ws.Range("H8").Value = "6:0"
$ws.Range("L8:M8").Value = "0"
$ws.Range("N8:O8").NumberFormat = "@"
$ws.Range("N8").Value = "100.00"
$ws.Range("P8").Value = "100.0000"
$ws.Range("Q8").NumberFormat = "@"
$ws.Range("Q8").Value = "1:0"

$ws.Rows(8).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 4. Row 9 - the new "total transactions" cell (merged P9:Q9, inherited from
#    the insert already); give it its value and correct row height.
# ---------------------------------------------------------------------------
$ws.Range("P9:Q9").Merge()
$ws.Range("P9").Value = 200
$ws.Rows(9).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 5. Row 10 (formerly row 9) - refresh the timestamp text to the new time.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Wednesday, 18 June, 2025 12:27 AM"
